$d = $word.ActiveDocument

# The "Module structure" list contains one paragraph per file, e.g.
#   <tab>mod_tab_traits.R
# Per the commit "Remove mod_tab_traits.R from documentation", that whole
# paragraph (its tab + text + the paragraph mark that ends it) must be
# removed, leaving the surrounding paragraphs (mod_derived_traits.R / ...
# and mod_repeatability.R / ...) directly adjacent, exactly as they were
# before this paragraph existed.

$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*mod_tab_traits.R*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    # Delete the whole paragraph, including its trailing paragraph mark,
    # so the paragraph is fully removed rather than left blank.
    $target.Range.Delete()
}
